$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"
